$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 7).Value = 1.099153532933269
    $ws.Cells.Item($row, 8).Value = 0.998
}
